$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("notose_71", "Sorenson", "Yes", "wzcagccwyfvnyucucvdafcbnseccycqcifckqrncizvecaxgab"),
  @("nomut_71",  "Sorenson", "No",  "wzcagccwyfvnyucucvdafcbnseccycqcifckqrncizvecaxgab"),
  @("notose_72", "Sorenson", "Yes", "wzcagciqcuwcuvcicoazycqmccvbtyevrcxsrstiszvfaxxgab"),
  @("nomut_72",  "Sorenson", "No",  "wzcagclecuocqkcpconcyclrujvbtpcvrypsrmtmnzvfcaxgab"),
  @("notmut_73", "Sorenson", "No",  "wzcagclecuocnkcpconcyclruzvjtpcvrypsrmtmnzvfcaxgab"),
  @("notose_73", "Sorenson", "Yes", "wzcagciqiuwcuvticoazycqmccvbtyevrcxsrstiszvfaxxgab"),
  @("notose_74", "Sorenson", "Yes", "wzcagcneduwrbkoovoncyfuovovojpcyulpizytbazvfcaxgab"),
  @("notose_75", "Sorenson", "Yes", "wzcagcneduwrbkoovoncyfuovovojpcyulpizytbazvfcaxgab"),
  @("nomut_74",  "Sorenson", "No",  "wzcagcqeduwzbkdpcorcyfuoczwojrcvulplpmthczvfcaxgab"),
  @("nomut_75",  "Sorenson", "No",  "wzcagcqeduwzbkdpcorcyfuoczwojrcvulplpmthczvfclxgab"),
  @("nomut_76",  "Sorenson", "No",  "wzcagcqeduwcrkqpeorcefwoczvojrcvulplpmtpazvfcaxgab"),
  @("notose_76", "Sorenson", "Yes", "wzcagcneduwrbkoovoncyfuovpvojpcyulpizyhbazvfcaxgab"),
  @("notose_77", "Sorenson", "Yes", "wzcagcneduwrbkoovoncyfuovpvojpcyulpizyhbazvfcaxgab"),
  @("nomut_77",  "Sorenson", "No",  "wzcagcqeduwcrkqpcorcefwocrsojrcvulplpmtpazvfeaxgab"),
  @("notose_78", "Sorenson", "Yes", "wzcagcnydfwrjxmokoncafuovpqcjqtyulpjzhhwazvfcaxgab"),
  @("nomut_78",  "Sorenson", "No",  "wzcaxcwgcuciufjnypvzywslkcvblicvlfweiseptzvfcajgab"),
  @("nomut_79",  "Sorenson", "No",  "wzcagcnednwibkohvoncnfuovpycjpcyulpuzthfazvfcaxgab"),
  @("notose_79", "Sorenson", "Yes", "wzcagcwydfwrjxookoncafuovpqcjqtyulpjzhhwazvfcaxgab"),
  @("nomut_80",  "Sorenson", "No",  "wzcagcwydfwryxoskoncafooqpqcxqtyulpjzhhwazvfcaxgab"),
  @("notose_80", "Sorenson", "Yes", "wzcagcwydfwrvxookoncafuovpqcjqtyulnjzhxwazvfcaxgab")
)

$startRow = 142
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $data.Count - 1
$ws.Range("D$lastRow").Select() | Out-Null
